# Apply the edits described by the diff:
#  - Student date of birth: 1998-03-17 -> 1998-12-03
#  - Table A (courses) rows: titles and ECTS credit numbers rotate/change,
#    and the Total recalculates
#  - Contact email + position for the Student signature row

$d = $word.ActiveDocument

# --- 1. Date of birth (Table 1, row 2, col 4) ---
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 4).Range.Text = "1998-12-03"

# --- 2. Table A courses (Table 2, rows 4-6, col 3 = title, col 5 = ECTS) ---
$t2 = $d.Tables.Item(2)

$t2.Cell(4, 3).Range.Text = "Advanced data"
$t2.Cell(4, 5).Range.Text = "5"

$t2.Cell(5, 3).Range.Text = "Java Application"
$t2.Cell(5, 5).Range.Text = "4"

$t2.Cell(6, 3).Range.Text = "Scrum"
$t2.Cell(6, 5).Range.Text = "8"

# Total row (Table 2, row 8, col 5)
$t2.Cell(8, 5).Range.Text = "Total : 17"

# --- 3. Student signature row (Table 5, row 3): email + position ---
$t5 = $d.Tables.Item(5)
$t5.Cell(3, 3).Range.Text = "antoine@gmail.com"
$t5.Cell(3, 4).Range.Text = "Student"
